# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Updates the "Datos actualizados" timestamp string
#  - Refreshes case counts for several countries
#  - Fixes ordering of 8 country pairs (name swapped back into correct row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 19 de Junio de 2020 a las 13:10'
$ws.Range("B7").Value = 382143
$ws.Range("C7").Value = 1052
$ws.Range("D7").Value = 205454
$ws.Range("E7").Value = 164079
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 12610
$ws.Range("B13").Value = 200262
$ws.Range("C13").Value = 2615
$ws.Range("D13").Value = 159192
$ws.Range("E13").Value = 31678
$ws.Range("G13").Value = 120
$ws.Range("H13").Value = 9392
$ws.Range("B22").Value = 85462
$ws.Range("C22").Value = 1021
$ws.Range("D22").Value = 65409
$ws.Range("E22").Value = 19960
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 93
$ws.Range("A32").Value = 'Emiratos Arabes Unidos'
$ws.Range("B32").Value = 44145
$ws.Range("C32").Value = 393
$ws.Range("D32").Value = 30996
$ws.Range("E32").Value = 12849
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 300
$ws.Range("A33").Value = 'Indonesia'
$ws.Range("B33").Value = 43803
$ws.Range("C33").Value = 1041
$ws.Range("D33").Value = 17349
$ws.Range("E33").Value = 24081
$ws.Range("G33").Value = 34
$ws.Range("H33").Value = 2373
$ws.Range("B40").Value = 31217
$ws.Range("C40").Value = 17
$ws.Range("E40").Value = 361
$ws.Range("A47").Value = 'Rumania'
$ws.Range("B47").Value = 23400
$ws.Range("C47").Value = 320
$ws.Range("D47").Value = 16555
$ws.Range("E47").Value = 5361
$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 1484
$ws.Range("A48").Value = 'Panama'
$ws.Range("B48").Value = 23351
$ws.Range("D48").Value = 13782
$ws.Range("E48").Value = 9094
$ws.Range("H48").Value = 475
$ws.Range("A71").Value = 'Nepal'
$ws.Range("B71").Value = 8274
$ws.Range("C71").Value = 426
$ws.Range("D71").Value = 1402
$ws.Range("E71").Value = 6850
$ws.Range("H71").Value = 22
$ws.Range("A72").Value = 'Sudan'
$ws.Range("B72").Value = 8020
$ws.Range("D72").Value = 2966
$ws.Range("E72").Value = 4567
$ws.Range("H72").Value = 487
$ws.Range("B77").Value = 5639
$ws.Range("C77").Value = 164
$ws.Range("D77").Value = 3788
$ws.Range("E77").Value = 1772
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 79
$ws.Range("B78").Value = 5477
$ws.Range("C78").Value = 194
$ws.Range("D78").Value = 719
$ws.Range("E78").Value = 4636
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 122
$ws.Range("A119").Value = 'Madagascar'
$ws.Range("B119").Value = 1443
$ws.Range("C119").Value = 40
$ws.Range("D119").Value = 498
$ws.Range("E119").Value = 932
$ws.Range("H119").Value = 13
$ws.Range("A120").Value = 'Zambia'
$ws.Range("B120").Value = 1416
$ws.Range("D120").Value = 1144
$ws.Range("E120").Value = 261
$ws.Range("H120").Value = 11
$ws.Range("B124").Value = 1128
$ws.Range("C124").Value = 3
$ws.Range("D124").Value = 1074
$ws.Range("E124").Value = 50
$ws.Range("B137").Value = 755
$ws.Range("C137").Value = 14
$ws.Range("E137").Value = 269
$ws.Range("D139").Value = 610
$ws.Range("E139").Value = 44
$ws.Range("D144").Value = 437
$ws.Range("E144").Value = 199
$ws.Range("B156").Value = 349
$ws.Range("C156").Value = 7
$ws.Range("E156").Value = 24
$ws.Range("A202").Value = 'Fiyi'
$ws.Range("A203").Value = 'Dominica'
$ws.Range("A206").Value = 'Groenlandia'
$ws.Range("A207").Value = 'Islas Malvinas'
$ws.Range("A210").Value = 'Montserrat'
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
$ws.Range("A211").Value = 'Seychelles'
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
$ws.Range("A213").Value = 'Islas Virgenes Britanicas'
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1
$ws.Range("A214").Value = 'Papua Nueva Guinea'
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
